# Insert one new data row at row 872, shifting existing rows 872:913 down to 873:914.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(872).Insert()

# Column A holds dates as plain text (e.g. "2026/02/24"); force text entry
# so Excel doesn't auto-convert the string into a date serial, then strip
# the temporary Text number-format so the cell carries no style index
# (matching the rest of the column).
$ws.Range("A872").NumberFormat = "@"
$ws.Range("A872").Value = "2026/02/24"
$ws.Range("A872").ClearFormats()

$ws.Range("B872").Value = "火"
$ws.Range("C872").Value = 10
$ws.Range("D872").Value = 201
